# BL-02 non-sync job retry diagram: fix AdminConnectionRetryInterceptor
# description (retry counter was showing "4th retry" text but graphic
# illustrated a 3rd failed retry followed by a successful 4th reconnect
# in one step too many) - drop the superfluous "3rd try...NG" callout
# + arrow from the retry-step group, nudge the surviving "2nd try" NG
# callout/arrow and the two red highlight boxes to the corrected spots,
# and reword the success caption to read "<halfwidth 3>回目のリトライで
# 再接続に成功！" instead of "４回目の...".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ------------------------------------------------------------------
# 1) The "1st/2nd/3rd try...NG" + dashed-arrow group (グループ化 148).
#    Re-position the 2nd try's callout box + its arrow, then delete the
#    3rd try's callout box + its arrow entirely.
# ------------------------------------------------------------------
$retryGroup = $s.Shapes.Item(3)

# Sanity anchor: this must be "グループ化 148" / id 149 containing the
# NG-retry callouts; GroupItems index order mirrors the original z-order
# (1=47,2=48,3=49,4=52,5=54,6=87 "1st try", 7=88 arrow#2, 8=89 "2nd try",
# 9=90 "3rd try", 10=91 arrow#3).
$arrow2 = $retryGroup.GroupItems.Item(7)
$arrow2.Top = 246.60157480314962

$callout2 = $retryGroup.GroupItems.Item(8)
$callout2.Left = 224.64063992125983
$callout2.Top = 224.79062992125984

# Ungroup so the 3rd-try callout + arrow can actually be removed (this
# engine's GroupShapes.Item(n).Delete() can't resolve shapes nested
# inside a group directly), then regroup the 8 survivors back together.
$flat = $retryGroup.Ungroup()
$flat.Item(10).Delete()   # 直線矢印コネクタ 90 (3rd try arrow)
$flat.Item(9).Delete()    # テキスト ボックス 89 ("３回目…NG")

$regrouped = $s.Shapes.Range(@(3, 4, 5, 6, 7, 8, 9, 10)).Group()
$regrouped.Name = "グループ化 148"

# ------------------------------------------------------------------
# 2) Red highlight rectangle over the retry counter - move it to line
#    up with the (now repositioned) "2nd try" callout.
# ------------------------------------------------------------------
$highlight1 = $s.Shapes.Item(6)
$highlight1.Left = 228.64552181102363
$highlight1.Top = 223.92188976377952

# ------------------------------------------------------------------
# 3) Success caption textbox: shrink it to fit the new wording and
#    reword "４回目のリトライで再接続に成功！" -> "<halfwidth 3>回目の
#    リトライで再接続に成功！" split across three runs (the halfwidth
#    "3" is tagged en-US, matching the source edit).
# ------------------------------------------------------------------
$caption = $s.Shapes.Item(7)
$caption.Width = 188.85142732283464

$capRange = $caption.TextFrame.TextRange
$capRange.Text = "3"
$run2 = $capRange.InsertAfter("回目")
$run3 = $run2.InsertAfter("のリトライで再接続に成功！")
$capRange.Characters(1, 1).LanguageID = "en-US"
